# Updated cryptos list (values per upstream diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.936.25"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "'1.815.87"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'309.50"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.4689"
$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("D8").Value = "'0.3694"
$ws.Range("E8").Value = "  -1.05%  "

$ws.Range("D9").Value = "'0.07371"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("D11").Value = "'20.39"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "'1.780.42"
$ws.Range("E12").Value = "  +3.32%  "

$ws.Range("D13").Value = "'5.382"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").Value = "'0.07084"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "'6.521"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").Value = "'91.95"
$ws.Range("E16").Value = "  +1.29%  "

$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").Value = "'14.73"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").Value = "'26.965.00"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").Value = "'5.330"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("D24").Value = "'2.042.41"
$ws.Range("E24").Value = "  +4.44%  "

$ws.Range("D25").Value = "'1.891"
$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("D26").Value = "'150.87"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("E27").Value = "  +1.80%  "

$ws.Range("D28").Value = "'18.38"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").Value = "'5.327"
$ws.Range("E29").Value = "  +0.59%  "

$ws.Range("D30").Value = "'116.19"
$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("D31").Value = "'0.08931"
$ws.Range("E31").Value = "  +0.39%  "

$ws.Range("D32").Value = "'0.7671"
$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("E33").Value = "  +0.69%  "

$ws.Range("D34").Value = "'4.509"
$ws.Range("E34").Value = "  +0.84%  "

$ws.Range("D35").Value = "'2.915"
$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").Value = "'1.086"
$ws.Range("E37").Value = "  -2.49%  "

$ws.Range("D38").Value = "'0.01963"
$ws.Range("E38").Value = "  +0.31%  "

$ws.Range("D39").Value = "'0.05294"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("D40").Value = "'2.973"
$ws.Range("E40").Value = "  +2.62%  "

$ws.Range("D41").Value = "'7.267"
$ws.Range("E41").Value = "  +0.60%  "

$ws.Range("D42").Value = "'0.5345"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "'2.328"
$ws.Range("E43").Value = "  -3.15%  "

$ws.Range("D44").Value = "'0.1653"
$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").Value = "'8.455"
$ws.Range("E45").Value = "  -1.25%  "

$ws.Range("D46").Value = "'0.4933"
$ws.Range("E46").Value = "  -2.37%  "

$ws.Range("D47").Value = "'10.45"
$ws.Range("E47").Value = "  +1.34%  "

$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").Value = "'1.671"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("D50").Value = "'103.25"
$ws.Range("E50").Value = "  -0.48%  "

$ws.Range("D51").Value = "'0.06304"
$ws.Range("E51").Value = "  -0.17%  "
